# Cleaned up code and commented + fixed picture packet parsing
#
# Rename the existing "Sheet1" to "Data packet" and add a new "Point packet"
# worksheet after it, describing the (smaller) point-packet structure.

$wb = $excel.ActiveWorkbook

# --- Rename the original sheet ---------------------------------------------
$dataSheet = $wb.Worksheets.Item(1)
$dataSheet.Name = "Data packet"

# --- Add the new "Point packet" sheet, right after "Data packet" -----------
$pointSheet = $wb.Worksheets.Add($null, $dataSheet)
$pointSheet.Name = "Point packet"

# Header row
$pointSheet.Range("A1").Value = "Data"
$pointSheet.Range("B1").Value = "Length"
$pointSheet.Range("C1").Value = "Position"
$pointSheet.Range("D1").Value = "Index"

# Data rows: label + length; Position/Index are computed with formulas,
# mirroring the layout used on the "Data packet" sheet.
$pointSheet.Range("A2").Value = "Start"
$pointSheet.Range("B2").Value = 2

$pointSheet.Range("A3").Value = "altAlt"
$pointSheet.Range("B3").Value = 4

$pointSheet.Range("A4").Value = "lat"
$pointSheet.Range("B4").Value = 4

$pointSheet.Range("A5").Value = "long"
$pointSheet.Range("B5").Value = 4

$pointSheet.Range("A6").Value = "altGPS"
$pointSheet.Range("B6").Value = 4

$pointSheet.Range("A7").Value = "heading"
$pointSheet.Range("B7").Value = 4

$pointSheet.Range("A8").Value = "End"
$pointSheet.Range("B8").Value = 2

# Running-total ("Position") and per-row "Index" formulas.
$pointSheet.Range("C2").Formula = "=SUM(B`$2:B2)"
$pointSheet.Range("D2").Formula = "=C2-B2"

$pointSheet.Range("C3").Formula = "=SUM(B`$2:B3)"
$pointSheet.Range("D3").Formula = "=C3-B3"

$pointSheet.Range("C4").Formula = "=SUM(B`$2:B4)"
$pointSheet.Range("D4").Formula = "=C4-B4"

$pointSheet.Range("C5").Formula = "=SUM(B`$2:B5)"
$pointSheet.Range("D5").Formula = "=C5-B5"

$pointSheet.Range("C6").Formula = "=SUM(B`$2:B6)"
$pointSheet.Range("D6").Formula = "=C6-B6"

$pointSheet.Range("C7").Formula = "=SUM(B`$2:B7)"
$pointSheet.Range("D7").Formula = "=C7-B7"

$pointSheet.Range("C8").Formula = "=SUM(B`$2:B8)"
$pointSheet.Range("D8").Formula = "=C8-B8"

# Leave selection on B2 and make "Point packet" the active/visible tab,
# matching the saved view state of the workbook.
$pointSheet.Range("B2").Select() | Out-Null
$pointSheet.Activate() | Out-Null
